$wb = $excel.ActiveWorkbook

# Add the new "Aggregate" worksheet after the last existing sheet (2010)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Aggregate"

# Header row
$ws.Range("A1").Value = "FIPS"
$ws.Range("B1").Value = "State"
$ws.Range("C1").Value = "EVER"
$ws.Range("D1").Value = "SAMPLE"
$ws.Range("E1").Value = "PRV per 100"

# National / by-state aggregate count data (EVER, SAMPLE, PRV per 100)
$ws.Range("A2").Value = 1.0
$ws.Range("B2").Value = "Alabama"
$ws.Range("C2").Value = 157450.00456785876
$ws.Range("D2").Value = 1091026.5426309681
$ws.Range("E2").Value = 14.431363345952535
$ws.Range("A3").Value = 4.0
$ws.Range("B3").Value = "Arizona"
$ws.Range("C3").Value = 866487.1122159892
$ws.Range("D3").Value = 6631643.360124848
$ws.Range("E3").Value = 13.065948591657628
$ws.Range("A4").Value = 6.0
$ws.Range("B4").Value = "California"
$ws.Range("C4").Value = 4513981.265956005
$ws.Range("D4").Value = 36944762.16889343
$ws.Range("E4").Value = 12.218190078800035
$ws.Range("A5").Value = 9.0
$ws.Range("B5").Value = "Connecticut"
$ws.Range("C5").Value = 515371.889014145
$ws.Range("D5").Value = 3216911.134948612
$ws.Range("E5").Value = 16.02070642906888
$ws.Range("A6").Value = 10.0
$ws.Range("B6").Value = "Delaware"
$ws.Range("C6").Value = 37621.98624757164
$ws.Range("D6").Value = 206483.30947034978
$ws.Range("E6").Value = 18.220352213491626
$ws.Range("A7").Value = 11.0
$ws.Range("B7").Value = "District of Columbia"
$ws.Range("C7").Value = 109403.39658694674
$ws.Range("D7").Value = 550985.4432505644
$ws.Range("E7").Value = 19.855950447895008
$ws.Range("A8").Value = 13.0
$ws.Range("B8").Value = "Geogia"
$ws.Range("C8").Value = 1847944.0559254412
$ws.Range("D8").Value = 12211231.561340636
$ws.Range("E8").Value = 15.133150547859733
$ws.Range("A9").Value = 16.0
$ws.Range("B9").Value = "Idaho"
$ws.Range("C9").Value = 104996.84186401556
$ws.Range("D9").Value = 1168919.0900528755
$ws.Range("E9").Value = 8.982387468688366
$ws.Range("A10").Value = 17.0
$ws.Range("B10").Value = "Illinois"
$ws.Range("C10").Value = 1580895.706567785
$ws.Range("D10").Value = 12758370.67075068
$ws.Range("E10").Value = 12.391047002514844
$ws.Range("A11").Value = 18.0
$ws.Range("B11").Value = "Indiana"
$ws.Range("C11").Value = 1006366.2885382811
$ws.Range("D11").Value = 7837909.728326183
$ws.Range("E11").Value = 12.839727981317218
$ws.Range("A12").Value = 19.0
$ws.Range("B12").Value = "Iowa"
$ws.Range("C12").Value = 287608.6462125769
$ws.Range("D12").Value = 3410013.627202459
$ws.Range("E12").Value = 8.434237444632387
$ws.Range("A13").Value = 20.0
$ws.Range("B13").Value = "Kansas"
$ws.Range("C13").Value = 396146.7413020134
$ws.Range("D13").Value = 3428397.615133414
$ws.Range("E13").Value = 11.554865735332672
$ws.Range("A14").Value = 21.0
$ws.Range("B14").Value = "Kentucky"
$ws.Range("C14").Value = 413087.5759674031
$ws.Range("D14").Value = 2957081.5267942287
$ws.Range("E14").Value = 13.969434803349209
$ws.Range("A15").Value = 22.0
$ws.Range("B15").Value = "Louisiana"
$ws.Range("C15").Value = 282442.7614159923
$ws.Range("D15").Value = 2178496.230470162
$ws.Range("E15").Value = 12.96503328605878
$ws.Range("A16").Value = 23.0
$ws.Range("B16").Value = "Maine"
$ws.Range("C16").Value = 109119.52465686374
$ws.Range("D16").Value = 825220.5687211744
$ws.Range("E16").Value = 13.223073780863679
$ws.Range("A17").Value = 24.0
$ws.Range("B17").Value = "Maryland"
$ws.Range("C17").Value = 1002975.984780239
$ws.Range("D17").Value = 6754689.3152769
$ws.Range("E17").Value = 14.848587965575783
$ws.Range("A18").Value = 26.0
$ws.Range("B18").Value = "Michigan"
$ws.Range("C18").Value = 1625547.03354681
$ws.Range("D18").Value = 11990510.426796142
$ws.Range("E18").Value = 13.55694608224577
$ws.Range("A19").Value = 27.0
$ws.Range("B19").Value = "Minnesota"
$ws.Range("C19").Value = 115805.68122609484
$ws.Range("D19").Value = 1222301.5589915775
$ws.Range("E19").Value = 9.474395281115143
$ws.Range("A20").Value = 28.0
$ws.Range("B20").Value = "Mississippi"
$ws.Range("C20").Value = 515551.3439893284
$ws.Range("D20").Value = 3630502.6189303007
$ws.Range("E20").Value = 14.20055011945513
$ws.Range("A21").Value = 29.0
$ws.Range("B21").Value = "Missouri"
$ws.Range("C21").Value = 766863.5236230671
$ws.Range("D21").Value = 5518463.724450555
$ws.Range("E21").Value = 13.896322634600988
$ws.Range("A22").Value = 30.0
$ws.Range("B22").Value = "Montana"
$ws.Range("C22").Value = 102944.19611154444
$ws.Range("D22").Value = 1059004.3401108915
$ws.Range("E22").Value = 9.720847423605917
$ws.Range("A23").Value = 31.0
$ws.Range("B23").Value = "Nebraska"
$ws.Range("C23").Value = 204904.80981110327
$ws.Range("D23").Value = 2201247.7197266174
$ws.Range("E23").Value = 9.308575676187472
$ws.Range("A24").Value = 32.0
$ws.Range("B24").Value = "Nevada"
$ws.Range("C24").Value = 351541.34141214454
$ws.Range("D24").Value = 3232267.2275993405
$ws.Range("E24").Value = 10.875998692510342
$ws.Range("A25").Value = 33.0
$ws.Range("B25").Value = "New Hampshire"
$ws.Range("C25").Value = 107549.14772854405
$ws.Range("D25").Value = 886427.4615532869
$ws.Range("E25").Value = 12.132876337121333
$ws.Range("A26").Value = 34.0
$ws.Range("B26").Value = "New Jersey"
$ws.Range("C26").Value = 1168380.2370796714
$ws.Range("D26").Value = 8196056.390538524
$ws.Range("E26").Value = 14.255395295089016
$ws.Range("A27").Value = 35.0
$ws.Range("B27").Value = "New Mexico"
$ws.Range("C27").Value = 179120.75376638884
$ws.Range("D27").Value = 1497759.9561924192
$ws.Range("E27").Value = 11.959243070014148
$ws.Range("A28").Value = 36.0
$ws.Range("B28").Value = "New York"
$ws.Range("C28").Value = 2781426.375356151
$ws.Range("D28").Value = 17587680.554876994
$ws.Range("E28").Value = 15.814628692382476
$ws.Range("A29").Value = 38.0
$ws.Range("B29").Value = "North Dakota"
$ws.Range("C29").Value = 37699.56334368641
$ws.Range("D29").Value = 425778.2815509131
$ws.Range("E29").Value = 8.854271102406717
$ws.Range("A30").Value = 39.0
$ws.Range("B30").Value = "Ohio"
$ws.Range("C30").Value = 997197.7980750023
$ws.Range("D30").Value = 8133869.580413682
$ws.Range("E30").Value = 12.259820350159654
$ws.Range("A31").Value = 40.0
$ws.Range("B31").Value = "Oklahoma"
$ws.Range("C31").Value = 487287.2730343682
$ws.Range("D31").Value = 3491912.5849495833
$ws.Range("E31").Value = 13.954738590382087
$ws.Range("A32").Value = 41.0
$ws.Range("B32").Value = "Oregon"
$ws.Range("C32").Value = 281481.12543031445
$ws.Range("D32").Value = 2525767.163670863
$ws.Range("E32").Value = 11.14438137762546
$ws.Range("A33").Value = 42.0
$ws.Range("B33").Value = "Pennsylvania"
$ws.Range("C33").Value = 1905109.2449985659
$ws.Range("D33").Value = 13667687.409110986
$ws.Range("E33").Value = 13.938782677518697
$ws.Range("A34").Value = 44.0
$ws.Range("B34").Value = "Rhode Island"
$ws.Range("C34").Value = 145915.22512547293
$ws.Range("D34").Value = 907042.8718127855
$ws.Range("E34").Value = 16.086916027888698
$ws.Range("A35").Value = 48.0
$ws.Range("B35").Value = "Texas"
$ws.Range("C35").Value = 3420043.6562700826
$ws.Range("D35").Value = 26030067.579157054
$ws.Range("E35").Value = 13.13881973556073
$ws.Range("A36").Value = 49.0
$ws.Range("B36").Value = "Utah"
$ws.Range("C36").Value = 404738.42055012786
$ws.Range("D36").Value = 3963226.5196104664
$ws.Range("E36").Value = 10.21234639371328
$ws.Range("A37").Value = 50.0
$ws.Range("B37").Value = "Vermont"
$ws.Range("C37").Value = 89457.40950244437
$ws.Range("D37").Value = 646238.5548435429
$ws.Range("E37").Value = 13.842784345186956
$ws.Range("A38").Value = 51.0
$ws.Range("B38").Value = "Virginia"
$ws.Range("C38").Value = 979105.1491755865
$ws.Range("D38").Value = 7200938.098480222
$ws.Range("E38").Value = 13.59691106610442
$ws.Range("A39").Value = 53.0
$ws.Range("B39").Value = "Washington"
$ws.Range("C39").Value = 330137.8695229291
$ws.Range("D39").Value = 3063863.432174641
$ws.Range("E39").Value = 10.775214915131084
$ws.Range("A40").Value = 54.0
$ws.Range("B40").Value = "West Virginia"
$ws.Range("C40").Value = 238133.07028626621
$ws.Range("D40").Value = 1877224.0092325318
$ws.Range("E40").Value = 12.685383796237643
$ws.Range("A41").Value = 55.0
$ws.Range("B41").Value = "Wisconsin"
$ws.Range("C41").Value = 410614.6794413228
$ws.Range("D41").Value = 3890402.9564449224
$ws.Range("E41").Value = 10.554553963647647
$ws.Range("A42").Value = 56.0
$ws.Range("B42").Value = "Wyoming"
$ws.Range("C42").Value = 34970.686195598726
$ws.Range("D42").Value = 367927.95107951324
$ws.Range("E42").Value = 9.504764748911718

# Restore the originally active sheet/tab selection
$wb.Worksheets.Item(1).Activate()
